$d = $word.ActiveDocument

# Locate the paragraph that contains the M2Doc field
# " m:null.setWidth(null) " (begin fldChar ... end fldChar).
$field = $d.Fields.Item(1)
$codeStart = $field.Code.Start

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($codeStart -ge $p.Range.Start -and $codeStart -lt $p.Range.End) {
        $targetParagraph = $p
    }
}

# Rebuild the paragraph as plain text runs "{m:null.setWidth(null)}" instead
# of a real Word field, keeping the orange coloring on the
# "null.setWidth(null" portion exactly as before, and preserving the
# paragraph's identity (paraId/textId/rsids).
$flatOpc = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="56541801" w14:textId="57C28A77" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>null</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.setWidth(</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>null</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetParagraph.Range.InsertXML($flatOpc)
